$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a throwaway sheet after the last sheet so it consumes the next internal
# sheetId (4), leaving sheetId 5 for the real new sheet we create below (to
# match the target workbook.xml: Gesamtübersicht=2, Monatsübersicht=5, Flags=3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$throwaway = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Duplicate "Gesamtübersicht" right after itself - this gives us an exact
# copy of all formatting/formulas/merged cells, placed directly before "Flags".
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "Monatsübersicht"

# Remove the throwaway sheet now that the real copy already grabbed sheetId 5.
$wb.Worksheets.Item("Sheet1").Delete()

# Restore "Gesamtübersicht" as the active/selected tab.
$wb.Worksheets.Item(1).Select()

# Register the hidden filter-database defined name for the new sheet, just
# like the one that already exists for "Gesamtübersicht".
$monats = $wb.Worksheets.Item("Monatsübersicht")
$monats.Names.Add("_xlnm._FilterDatabase", "=Monatsübersicht!`$B`$3:`$I`$3") | Out-Null
$filterName = $wb.Names.Item("Monatsübersicht!_FilterDatabase")
$filterName.Visible = $false
